$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O2").Value = "Lucy Padilla"
$ws.Range("S2").Value = "7900 Colony Cir S"
$ws.Range("T2").Value = "Apt 305"
$ws.Range("V2").Value = "Tamarac"
$ws.Range("W2").Value = "FL"
$ws.Range("Y2").Value = "'33321"
$ws.Range("AD2").Value = '"Gold Filled Beads"'
$ws.Range("AF2").Value = "'16.00"
